$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F8").Value = -3
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = -4
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 5
$ws.Range("F30").Value = -2
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = 3
$ws.Range("F47").Value = 1
$ws.Range("F49").Value = 0
$ws.Range("F61").Value = -5
$ws.Range("F63").Value = 3
